# Inserts two new price records (rows 550 and 551) for date 2022-02-18
# (serial 44610) for Brócoli at Vega Central Mapocho de Santiago,
# shifting all subsequent rows down by two (old row 550 -> new row 552,
# ..., old row 581 -> new row 583), matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 550; everything from the old row 550
# onward shifts down to make room (Excel copies formatting, e.g. the
# date NumberFormat on column D, from the row above automatically).
$ws.Rows("550:551").Insert()

# New row 550 - "Primera" grade entry for 2022-02-18
$ws.Range("A550").Value2 = 9
$ws.Range("B550").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C550").Value = "Metropolitana"
$ws.Range("D550").Value2 = 44610
$ws.Range("E550").Value2 = 13
$ws.Range("F550").Value2 = 100112023
$ws.Range("G550").Value = "Brócoli"
$ws.Range("H550").Value = "Sin especificar"
$ws.Range("I550").Value = "Primera"
$ws.Range("J550").Value2 = 4300
$ws.Range("K550").Value2 = 900
$ws.Range("L550").Value2 = 1000
$ws.Range("M550").Value2 = 950
$ws.Range("N550").Value = "`$/unidad"
$ws.Range("O550").Value = "Región Metropolitana"
$ws.Range("P550").Value2 = 950
$ws.Range("Q550").Value2 = 1
$ws.Range("R550").Value = "Hortaliza"

# New row 551 - "Segunda" grade entry for 2022-02-18
$ws.Range("A551").Value2 = 9
$ws.Range("B551").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C551").Value = "Metropolitana"
$ws.Range("D551").Value2 = 44610
$ws.Range("E551").Value2 = 13
$ws.Range("F551").Value2 = 100112023
$ws.Range("G551").Value = "Brócoli"
$ws.Range("H551").Value = "Sin especificar"
$ws.Range("I551").Value = "Segunda"
$ws.Range("J551").Value2 = 1960
$ws.Range("K551").Value2 = 700
$ws.Range("L551").Value2 = 700
$ws.Range("M551").Value2 = 700
$ws.Range("N551").Value = "`$/unidad"
$ws.Range("O551").Value = "Región Metropolitana"
$ws.Range("P551").Value2 = 700
$ws.Range("Q551").Value2 = 1
$ws.Range("R551").Value = "Hortaliza"
